# Add I0 (column I) and IF (column J) data to the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells before setting their text, so I1/J1 end up with the same
# bold/centered/bordered style used by the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (I0) and column J (IF), rows 2-72
$iValues = @(5,7,9,7,6,7,10,7,6,12,7,8,7,7,8,6,6,9,9,9,9,9,8,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,8,7,8,8,8,7,3,7,5,6,4,4,4)
$jValues = @(5,7,9,7,7,7,10,7,6,12,7,8,7,7,8,6,7,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,8,8,8,8,8,7,3,7,5,6,4,4,4)

for ($r = 2; $r -le 72; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
